$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 862
$ws.Range("F8").Value = 48
$ws.Range("F11").Value = 295
$ws.Range("F12").Value = 296
$ws.Range("F14").Value = 239
$ws.Range("F17").Value = 6627
$ws.Range("F18").Value = 64
$ws.Range("F19").Value = 72
$ws.Range("F21").Value = 7577
$ws.Range("F22").Value = 38
$ws.Range("F24").Value = 3393
$ws.Range("F25").Value = 28
$ws.Range("F26").Value = 1797
$ws.Range("F27").Value = 893
$ws.Range("F28").Value = 4515
$ws.Range("F29").Value = 125
$ws.Range("F30").Value = 351
$ws.Range("F34").Value = 1672
$ws.Range("F36").Value = 164
$ws.Range("F39").Value = 1203
$ws.Range("F40").Value = 1758
$ws.Range("F41").Value = 2135

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 69
$ws.Range("F5").Value = 6

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 243
$ws.Range("F3").Value = 1225

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 243
$ws.Range("F4").Value = 1225
$ws.Range("F9").Value = 862
$ws.Range("F10").Value = 48
$ws.Range("F13").Value = 295
$ws.Range("F14").Value = 297
$ws.Range("F15").Value = 69
$ws.Range("F17").Value = 239
$ws.Range("F20").Value = 6627
$ws.Range("F21").Value = 64
$ws.Range("F22").Value = 72
$ws.Range("F24").Value = 7577
$ws.Range("F25").Value = 38
$ws.Range("F27").Value = 3393
$ws.Range("F28").Value = 28
$ws.Range("F29").Value = 1797
$ws.Range("F30").Value = 893
$ws.Range("F31").Value = 4515
$ws.Range("F32").Value = 125
$ws.Range("F33").Value = 351
$ws.Range("F38").Value = 1672
$ws.Range("F40").Value = 164
$ws.Range("F43").Value = 6
$ws.Range("F44").Value = 1203
$ws.Range("F45").Value = 1758
$ws.Range("F47").Value = 2135
